$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row suffixes to reflect the respective input file names:
#   "<name>_old" -> "<name>_FV2310"
#   "<name>_new" -> "<name>_FV2404"
# The "diff" header (column K / 11) stays as-is.
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2
    if ($text -ne $null) {
        if ($text.EndsWith("_old")) {
            $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2310"
        } elseif ($text.EndsWith("_new")) {
            $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2404"
        }
    }
}

# Add a table over A1:U66 with autofilter
$rng = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze panes at row 2 (freeze the header row)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
